{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Paragraph: 'Link \"Consultar novo veiculo\" ...' -> make all text red (FF0000),\n//    including the paragraph mark itself.\nconst p1 = paragraphs.items[1];\nif (!p1.text.includes(\"Consultar novo\")) {\n  throw new Error(\"Paragraph 1 text mismatch\");\n}\np1.font.color = \"#FF0000\";\n\n// 2) Paragraph: 'Se o servidor visualizar uma reserva ...' -> make all text red (FF0000),\n//    including the paragraph mark itself.\nconst p2 = paragraphs.items[2];\nif (!p2.text.includes(\"servidor visualizar\")) {\n  throw new Error(\"Paragraph 2 text mismatch\");\n}\np2.font.color = \"#FF0000\";\n\n// 3) Paragraph: 'No alterar reserva, nao e possivel alterar a data ...' -> append a\n//    new trailing run containing just a single space.\nconst p3 = paragraphs.items[7];\nif (!p3.text.includes(\"No alterar reserva\")) {\n  throw new Error(\"Paragraph 7 text mismatch\");\n}\np3.insertText(\" \", Word.InsertLocation.end);\n\n// 4) Paragraph: 'Quando vai alterar a senha, se o campo senha atual estiver errada ...'\n//    -> make all text red (FF0000), including the paragraph mark itself.\nconst p4 = paragraphs.items[9];\nif (!p4.text.includes(\"senha atual estiver errada\")) {\n  throw new Error(\"Paragraph 9 text mismatch\");\n}\np4.font.color = \"#FF0000\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Paragraph: 'Link \"Consultar novo veiculo\" ...' -> make all text red (FF0000),\n#    including the paragraph mark itself.\n$p1 = $d.Paragraphs(2)\nif ($p1.Range.Text -notlike \"*Consultar novo*\") { throw \"Paragraph 2 text mismatch\" }\n$p1.Range.Font.Color = 255\n\n# 2) Paragraph: 'Se o servidor visualizar uma reserva ...' -> make all text red (FF0000),\n#    including the paragraph mark itself.\n$p2 = $d.Paragraphs(3)\nif ($p2.Range.Text -notlike \"*servidor visualizar*\") { throw \"Paragraph 3 text mismatch\" }\n$p2.Range.Font.Color = 255\n\n# 3) Paragraph: 'No alterar reserva, nao e possivel alterar a data ...' -> append a\n#    new trailing run containing just a single space.\n$p3 = $d.Paragraphs(8)\nif ($p3.Range.Text -notlike \"*No alterar reserva*\") { throw \"Paragraph 8 text mismatch\" }\n$r3 = $p3.Range\n$r3.SetRange($r3.End - 1, $r3.End - 1)\n$r3.InsertAfter(\" \")\n\n# 4) Paragraph: 'Quando vai alterar a senha, se o campo senha atual estiver errada ...'\n#    -> make all text red (FF0000), including the paragraph mark itself.\n$p4 = $d.Paragraphs(10)\nif ($p4.Range.Text -notlike \"*senha atual estiver errada*\") { throw \"Paragraph 10 text mismatch\" }\n$p4.Range.Font.Color = 255\n"}
